$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1942446043165468
$ws.Range("C2").Value = 0.5863309352517986
$ws.Range("J2").Value = 0.01079136690647482
$ws.Range("P2").Value = 0.1330935251798561
$ws.Range("S2").Value = 0.07553956834532374
$ws.Range("B3").Value = 0.005917159763313609
$ws.Range("C3").Value = 0.02366863905325444
$ws.Range("J3").Value = 0.02958579881656805
$ws.Range("P3").Value = 0.7869822485207101
$ws.Range("S3").Value = 0.1538461538461539
$ws.Range("P4").Value = 0.7906976744186046
$ws.Range("S4").Value = 0.2093023255813954
$ws.Range("B6").Value = 0.04545454545454546
$ws.Range("D6").Value = 0.01363636363636364
$ws.Range("E6").Value = 0.004545454545454545
$ws.Range("F6").Value = 0.04545454545454546
$ws.Range("J6").Value = 0.2863636363636364
$ws.Range("O6").Value = 0.04090909090909091
$ws.Range("Q6").Value = 0.1545454545454545
$ws.Range("R6").Value = 0.06818181818181818
$ws.Range("S6").Value = 0.3409090909090909
$ws.Range("B7").Value = 0.07575757575757576
$ws.Range("D7").Value = 0.03535353535353535
$ws.Range("F7").Value = 0.04545454545454546
$ws.Range("J7").Value = 0.1515151515151515
$ws.Range("O7").Value = 0.01515151515151515
$ws.Range("Q7").Value = 0.1868686868686869
$ws.Range("R7").Value = 0.08080808080808081
$ws.Range("S7").Value = 0.4090909090909091
$ws.Range("B8").Value = 0.07932692307692307
$ws.Range("D8").Value = 0.01682692307692308
$ws.Range("E8").Value = 0.002403846153846154
$ws.Range("F8").Value = 0.0625
$ws.Range("J8").Value = 0.1201923076923077
$ws.Range("O8").Value = 0.007211538461538462
$ws.Range("Q8").Value = 0.21875
$ws.Range("R8").Value = 0.0673076923076923
$ws.Range("S8").Value = 0.4254807692307692
$ws.Range("B9").Value = 0.106280193236715
$ws.Range("D9").Value = 0.02415458937198068
$ws.Range("F9").Value = 0.08695652173913043
$ws.Range("J9").Value = 0.09178743961352658
$ws.Range("O9").Value = 0.00966183574879227
$ws.Range("Q9").Value = 0.1980676328502415
$ws.Range("R9").Value = 0.06280193236714976
$ws.Range("S9").Value = 0.4202898550724637
$ws.Range("B10").Value = 0.1086286594761171
$ws.Range("D10").Value = 0.01926040061633282
$ws.Range("E10").Value = 0.0007704160246533128
$ws.Range("F10").Value = 0.06779661016949153
$ws.Range("J10").Value = 0.1271186440677966
$ws.Range("O10").Value = 0.01540832049306626
$ws.Range("Q10").Value = 0.2349768875192604
$ws.Range("R10").Value = 0.06317411402157165
$ws.Range("S10").Value = 0.3628659476117103
$ws.Range("G11").Value = 0.1501597444089457
$ws.Range("J11").Value = 0.1054313099041534
$ws.Range("K11").Value = 0.194888178913738
$ws.Range("L11").Value = 0.5335463258785943
$ws.Range("S11").Value = 0.01597444089456869
$ws.Range("G12").Value = 0.7078651685393258
$ws.Range("J12").Value = 0.2247191011235955
$ws.Range("K12").Value = 0.01123595505617977
$ws.Range("L12").Value = 0.02808988764044944
$ws.Range("S12").Value = 0.02808988764044944
$ws.Range("G13").Value = 0.7045454545454546
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.04545454545454546
$ws.Range("F15").Value = 0.01932367149758454
$ws.Range("H15").Value = 0.1256038647342995
$ws.Range("I15").Value = 0.07729468599033816
$ws.Range("J15").Value = 0.3671497584541063
$ws.Range("K15").Value = 0.07246376811594203
$ws.Range("M15").Value = 0.00966183574879227
$ws.Range("O15").Value = 0.0821256038647343
$ws.Range("S15").Value = 0.2463768115942029
$ws.Range("F16").Value = 0.01492537313432836
$ws.Range("H16").Value = 0.1741293532338309
$ws.Range("I16").Value = 0.08955223880597014
$ws.Range("J16").Value = 0.4527363184079602
$ws.Range("K16").Value = 0.109452736318408
$ws.Range("M16").Value = 0.01492537313432836
$ws.Range("O16").Value = 0.05970149253731343
$ws.Range("S16").Value = 0.0845771144278607
$ws.Range("F17").Value = 0.009881422924901186
$ws.Range("H17").Value = 0.1600790513833992
$ws.Range("I17").Value = 0.1027667984189723
$ws.Range("J17").Value = 0.4426877470355731
$ws.Range("K17").Value = 0.07509881422924901
$ws.Range("M17").Value = 0.02371541501976284
$ws.Range("O17").Value = 0.06521739130434782
$ws.Range("S17").Value = 0.1205533596837945
$ws.Range("F18").Value = 0.03947368421052631
$ws.Range("H18").Value = 0.2368421052631579
$ws.Range("I18").Value = 0.1052631578947368
$ws.Range("J18").Value = 0.3947368421052632
$ws.Range("K18").Value = 0.09868421052631579
$ws.Range("M18").Value = 0.01973684210526316
$ws.Range("O18").Value = 0.03947368421052631
$ws.Range("S18").Value = 0.06578947368421052
$ws.Range("F19").Value = 0.01868399675060926
$ws.Range("H19").Value = 0.1974004874086109
$ws.Range("I19").Value = 0.08529650690495533
$ws.Range("J19").Value = 0.3655564581640942
$ws.Range("K19").Value = 0.1251015434606012
$ws.Range("M19").Value = 0.02030869212022746
$ws.Range("N19").Value = 0.0008123476848090983
$ws.Range("O19").Value = 0.06417546709991877
$ws.Range("S19").Value = 0.1226645004061738
